# Fix minutes/seconds formatting in the "Общее время" (total time / haul) column:
# values like "N ч. M мин. S сек." should have M and S zero-padded to two digits
# (hours are left as-is). E.g. "55 ч. 45 мин. 0 сек." -> "55 ч. 45 мин. 00 сек."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

# Column 9 = "I" = "Общее время" (total time)
$col = 9

for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value2

    if ($val -match '^(\d+) ч\. (\d+) мин\. (\d+) сек\.$') {
        $hours = $matches[1]
        $minutes = $matches[2]
        $seconds = $matches[3]

        $minutesPadded = $minutes.PadLeft(2, '0')
        $secondsPadded = $seconds.PadLeft(2, '0')

        $newVal = "$hours ч. $minutesPadded мин. $secondsPadded сек."

        if ($newVal -ne $val) {
            $cell.Value2 = $newVal
        }
    }
}
